# Update the "Basic Commands" sheet: the sns.relplot() answers that used
# col="D" (double-quoted kwargs) are restated with col='D' (single-quoted
# kwargs), matching the rest of the sheet's quoting convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic Commands")

$ws.Range("B200").Value = "sns.relplot(data=df,x='A',y='B',hue='C',col='D')"
$ws.Range("B201").Value = "sns.relplot(data=df,x='A',y='B',hue='C',col='D',row='E')"
$ws.Range("B202").Value = "sns.relplot(data=df,x='A',y='B',hue='C',col='D',col_wrap=2)"
$ws.Range("B203").Value = "sns.relplot(data=df,x='A',y='B',hue='C',col='D',size='E',style='F')"
$ws.Range("B205").Value = "sns.relplot(data=df,x='A',y='B',hue='C',col='D',kind='line')"
$ws.Range("B206").Value = "sns.relplot(data=df,x='A',y='B',hue='C',col='D',row='E',kind='line')"
$ws.Range("B207").Value = "sns.relplot(data=df,x='A',y='B',hue='C',col='D',col_wrap=2,kind='line')"
$ws.Range("B208").Value = "sns.relplot(data=df,x='A',y='B',hue='C',col='D',size='E',style='F',kind='line')"

# Restore the author's final scroll position / active-cell selection.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 203
$win.ScrollColumn = 1
$ws.Range("B225").Select() | Out-Null

Write-Output "done"
